$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 188, shifting existing rows 188:206 down to 189:207
$ws.Rows.Item(188).Insert()

# Populate the newly inserted row 188 with the new record's data
$ws.Range("A188").Value = 11
$ws.Range("B188").Value = "Vega Monumental Concepción"
$ws.Range("C188").Value = "Bíobío"
$ws.Range("D188").Value = 45077
$ws.Range("E188").Value = 8
$ws.Range("F188").Value = "Fruta"
$ws.Range("G188").Value = 100102
$ws.Range("H188").Value = "Cítricos"
$ws.Range("I188").Value = 100102004
$ws.Range("J188").Value = "Mandarina"
$ws.Range("K188").Value = "Clementina"
$ws.Range("L188").Value = "Primera"
$ws.Range("M188").Value = 220
$ws.Range("N188").Value = 12000
$ws.Range("O188").Value = 13000
$ws.Range("P188").Value = 12455
$ws.Range("Q188").Value = "$/bandeja 18 kilos"
$ws.Range("R188").Value = "Provincia de Limarí"
$ws.Range("S188").Value = 692
$ws.Range("T188").Value = 18
